# Update vm_pu.xlsx bus voltage results for the 380 kV case (Case_4_48)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.044539588753954
$ws.Range("D2").Value = 1.050253013928409
$ws.Range("E2").Value = 1.048203787235186
$ws.Range("F2").Value = 1.057995675739932
$ws.Range("I2").Value = 1.039531495334753
$ws.Range("J2").Value = 1.049604130714789
$ws.Range("K2").Value = 1.053007805597117
$ws.Range("L2").Value = 1.050964286196767
$ws.Range("M2").Value = 1.060729119545004
$ws.Range("N2").Value = 1.020145131673962

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.046100966506257
$ws.Range("D3").Value = 1.051505264049348
$ws.Range("E3").Value = 1.049715425970514
$ws.Range("F3").Value = 1.05946515882987
$ws.Range("I3").Value = 1.039945920172026
$ws.Range("J3").Value = 1.050809779199649
$ws.Range("K3").Value = 1.054071331433928
$ws.Range("L3").Value = 1.052286112823378
$ws.Range("M3").Value = 1.06201088600449
$ws.Range("N3").Value = 1.020568724389414

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.047109321428878
$ws.Range("D4").Value = 1.052313610395426
$ws.Range("E4").Value = 1.050691973350411
$ws.Range("F4").Value = 1.060414398154055
$ws.Range("I4").Value = 1.040211526267952
$ws.Range("J4").Value = 1.051587598333353
$ws.Range("K4").Value = 1.054757007822976
$ws.Range("L4").Value = 1.0531393426588
$ws.Range("M4").Value = 1.062838149795622
$ws.Range("N4").Value = 1.020841451451857

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.04753277412096
$ws.Range("D5").Value = 1.052652980788161
$ws.Range("E5").Value = 1.051102143319062
$ws.Range("F5").Value = 1.060813079879312
$ws.Range("I5").Value = 1.040322578430404
$ws.Range("J5").Value = 1.05191404640866
$ws.Range("K5").Value = 1.055044674982014
$ws.Range("L5").Value = 1.053497550325911
$ws.Range("M5").Value = 1.063185430159601
$ws.Range("N5").Value = 1.020955781566558

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.047603847014526
$ws.Range("D6").Value = 1.052709935885507
$ws.Range("E6").Value = 1.051170991182768
$ws.Range("F6").Value = 1.060879998332781
$ws.Range("I6").Value = 1.040341188970159
$ws.Range("J6").Value = 1.051968826608854
$ws.Range("K6").Value = 1.055092941079095
$ws.Range("L6").Value = 1.053557666454705
$ws.Range("M6").Value = 1.063243710874946
$ws.Range("N6").Value = 1.020974959152141

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.04711498141885
$ws.Range("D7").Value = 1.052318146872097
$ws.Range("E7").Value = 1.050697455506933
$ws.Range("F7").Value = 1.060419726836076
$ws.Range("I7").Value = 1.040213012540363
$ws.Range("J7").Value = 1.051591962489105
$ws.Range("K7").Value = 1.054760853960834
$ws.Range("L7").Value = 1.053144130963428
$ws.Range("M7").Value = 1.062842792133689
$ws.Range("N7").Value = 1.020842980407648

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.045067675562221
$ws.Range("D8").Value = 1.050676624184095
$ws.Range("E8").Value = 1.048714985819223
$ws.Range("F8").Value = 1.058492633571165
$ws.Range("I8").Value = 1.039672082831048
$ws.Range("J8").Value = 1.050012069277105
$ws.Range("K8").Value = 1.053367749764966
$ws.Range("L8").Value = 1.051411439216589
$ws.Range("M8").Value = 1.061162742919995
$ws.Range("N8").Value = 1.020288571185299

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.041444605585405
$ws.Range("D9").Value = 1.047768862295522
$ws.Range("E9").Value = 1.045209100231216
$ws.Range("F9").Value = 1.055084122421749
$ws.Range("I9").Value = 1.038699199823196
$ws.Range("J9").Value = 1.047210025093734
$ws.Range("K9").Value = 1.050893520501972
$ws.Range("L9").Value = 1.048341926517304
$ws.Range("M9").Value = 1.058185676811884
$ws.Range("N9").Value = 1.019301053341025

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.039018220121151
$ws.Range("D10").Value = 1.045819710376649
$ws.Range("E10").Value = 1.042862877743305
$ws.Range("F10").Value = 1.05280271653045
$ws.Range("I10").Value = 1.038037183017026
$ws.Range("J10").Value = 1.045329378886133
$ws.Range("K10").Value = 1.049230587297231
$ws.Range("L10").Value = 1.046284133941128
$ws.Range("M10").Value = 1.056189335954531
$ws.Range("N10").Value = 1.018635430165579

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.037964820721998
$ws.Range("D11").Value = 1.044973084492009
$ws.Range("E11").Value = 1.041844688488678
$ws.Range("F11").Value = 1.051812578439176
$ws.Range("I11").Value = 1.037747295489037
$ws.Range("J11").Value = 1.044511941080483
$ws.Range("K11").Value = 1.048507238997456
$ws.Range("L11").Value = 1.045390259191386
$ws.Range("M11").Value = 1.055322035248509
$ws.Range("N11").Value = 1.018345444525067

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.037573114935938
$ws.Range("D12").Value = 1.044658206320643
$ws.Range("E12").Value = 1.041466137567472
$ws.Range("F12").Value = 1.051444445460449
$ws.Range("I12").Value = 1.037639129402256
$ws.Range("J12").Value = 1.044207832590384
$ws.Range("K12").Value = 1.048238053468194
$ws.Range("L12").Value = 1.045057798355792
$ws.Range("M12").Value = 1.054999440564348
$ws.Range("N12").Value = 1.018237462311262

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.037657156642264
$ws.Range("D13").Value = 1.044725767158705
$ws.Range("E13").Value = 1.041547354039639
$ws.Range("F13").Value = 1.05152342730522
$ws.Range("I13").Value = 1.037662353580947
$ws.Range("J13").Value = 1.044273086567146
$ws.Range("K13").Value = 1.048295817511389
$ws.Range("L13").Value = 1.045129132218577
$ws.Range("M13").Value = 1.055068658303768
$ws.Range("N13").Value = 1.018260637071503

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.037932450966184
$ws.Range("D14").Value = 1.044947064847126
$ws.Range("E14").Value = 1.041813404555231
$ws.Range("F14").Value = 1.051782155682688
$ws.Range("I14").Value = 1.037738364442566
$ws.Range("J14").Value = 1.044486813134832
$ws.Range("K14").Value = 1.048484998319392
$ws.Range("L14").Value = 1.04536278682493
$ws.Range("M14").Value = 1.055295378512096
$ws.Range("N14").Value = 1.018336524180578

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.038102012031634
$ws.Range("D15").Value = 1.045083359918438
$ws.Range("E15").Value = 1.041977280347335
$ws.Range("F15").Value = 1.051941519885602
$ws.Range("I15").Value = 1.037785132340911
$ws.Range("J15").Value = 1.044618433819067
$ws.Range("K15").Value = 1.048601492132972
$ws.Range("L15").Value = 1.045506691078782
$ws.Range("M15").Value = 1.055435009659318
$ws.Range("N15").Value = 1.018383245038524

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.039088071637302
$ws.Range("D16").Value = 1.045875842012158
$ws.Range("E16").Value = 1.042930403100934
$ws.Range("F16").Value = 1.052868379973469
$ws.Range("I16").Value = 1.038056353532878
$ws.Range("J16").Value = 1.04538356333404
$ws.Range("K16").Value = 1.049278523611234
$ws.Range("L16").Value = 1.046343396863993
$ws.Range("M16").Value = 1.05624683456203
$ws.Range("N16").Value = 1.018654638064296

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.039705853873689
$ws.Range("D17").Value = 1.046372234957257
$ws.Range("E17").Value = 1.043527659425084
$ws.Range("F17").Value = 1.053449158954435
$ws.Range("I17").Value = 1.038225616091442
$ws.Range("J17").Value = 1.045862671003018
$ws.Range("K17").Value = 1.049702321933796
$ws.Range("L17").Value = 1.046867473915764
$ws.Range("M17").Value = 1.056755295571904
$ws.Range("N17").Value = 1.018824400533433

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.040065930457327
$ws.Range("D18").Value = 1.046661519359866
$ws.Range("E18").Value = 1.043875811655792
$ws.Range("F18").Value = 1.053787698747906
$ws.Range("I18").Value = 1.038324032693606
$ws.Range("J18").Value = 1.046141827376647
$ws.Range("K18").Value = 1.049949199691011
$ws.Range("L18").Value = 1.047172886127648
$ws.Range("M18").Value = 1.057051595796093
$ws.Range("N18").Value = 1.018923249897567

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.040188662619115
$ws.Range("D19").Value = 1.046760115171957
$ws.Range("E19").Value = 1.043994486035524
$ws.Range("F19").Value = 1.053903095261773
$ws.Range("I19").Value = 1.038357537506541
$ws.Range("J19").Value = 1.046236962040618
$ws.Range("K19").Value = 1.050033325157027
$ws.Range("L19").Value = 1.047276977775437
$ws.Range("M19").Value = 1.057152579893611
$ws.Range("N19").Value = 1.018956926203197

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.039639599218208
$ws.Range("D20").Value = 1.046319002943346
$ws.Range("E20").Value = 1.043463602025012
$ws.Range("F20").Value = 1.053386869565389
$ws.Range("I20").Value = 1.038207488055519
$ws.Range("J20").Value = 1.045811298279743
$ws.Range("K20").Value = 1.049656885192639
$ws.Range("L20").Value = 1.046811273722789
$ws.Range("M20").Value = 1.056700771181148
$ws.Range("N20").Value = 1.018806204256478

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.0378513955148
$ws.Range("D21").Value = 1.044881909403863
$ws.Range("E21").Value = 1.041735069059873
$ws.Range("F21").Value = 1.051705976410333
$ws.Range("I21").Value = 1.037715994665699
$ws.Range("J21").Value = 1.044423889192697
$ws.Range("K21").Value = 1.048429303201541
$ws.Range("L21").Value = 1.045293993480082
$ws.Range("M21").Value = 1.055228627297495
$ws.Range("N21").Value = 1.018314184760296

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.03672460818288
$ws.Range("D22").Value = 1.043976012175838
$ws.Range("E22").Value = 1.040646240882834
$ws.Range("F22").Value = 1.050647092904772
$ws.Range("I22").Value = 1.03740414247952
$ws.Range("J22").Value = 1.043548813246836
$ws.Range("K22").Value = 1.047654566665618
$ws.Range("L22").Value = 1.044337491992467
$ws.Range("M22").Value = 1.054300478342009
$ws.Range("N22").Value = 1.01800327670258

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.03732217785497
$ws.Range("D23").Value = 1.044456470037392
$ws.Range("E23").Value = 1.041223645384201
$ws.Range("F23").Value = 1.051208623720339
$ws.Range("I23").Value = 1.037569730794493
$ws.Range("J23").Value = 1.044012971938334
$ws.Range("K23").Value = 1.048065547367863
$ws.Range("L23").Value = 1.044844794346861
$ws.Range("M23").Value = 1.054792752785818
$ws.Range("N23").Value = 1.018168243518011

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.039669537665088
$ws.Range("D24").Value = 1.046343056981632
$ws.Range("E24").Value = 1.043492547478289
$ws.Range("F24").Value = 1.053415016135197
$ws.Range("I24").Value = 1.038215680296868
$ws.Range("J24").Value = 1.045834512326142
$ws.Range("K24").Value = 1.049677417078665
$ws.Range("L24").Value = 1.04683666901475
$ws.Range("M24").Value = 1.056725409260798
$ws.Range("N24").Value = 1.018814426896784

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.042383150477735
$ws.Range("D25").Value = 1.048522432620823
$ws.Range("E25").Value = 1.046116995506354
$ws.Range("F25").Value = 1.055966864488165
$ws.Range("I25").Value = 1.038953066451098
$ws.Range("J25").Value = 1.047936609966844
$ws.Range("K25").Value = 1.051535508463988
$ws.Range("L25").Value = 1.049137450530468
$ws.Range("M25").Value = 1.058957335220895
$ws.Range("N25").Value = 1.019557620915198
